$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.358.48"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "1.573.07"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'211.99"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").Value = "'0.490"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D8").Value = "'44.60"
$ws.Range("E8").Value = "  -4.73%  "

$ws.Range("D9").Value = "'23.76"
$ws.Range("E9").Value = "  -1.49%  "

$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("E11").Value = "  -0.93%  "

$ws.Range("D12").Value = "'0.0897"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("D13").Value = "1.799.77"
$ws.Range("E13").Value = "  -0.31%  "

$ws.Range("D14").Value = "1.584.38"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("D16").Value = "'0.517"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("D17").Value = "28.366.46"
$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").Value = "'61.56"
$ws.Range("E18").Value = "  -1.32%  "

$ws.Range("D19").Value = "'230.41"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").Value = "'7.42"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "0.0₃0684"
$ws.Range("E21").Value = "  -1.60%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "'3.97"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("D25").Value = "'2.04"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").Value = "'151.52"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").Value = "'14.91"
$ws.Range("E27").Value = "  -0.81%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'6.37"
$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.104"
$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  +3.40%  "

$ws.Range("E32").Value = "  -3.80%  "

$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("E34").Value = "  -1.04%  "

$ws.Range("D35").Value = "1.386.02"
$ws.Range("E35").Value = "  -1.03%  "

$ws.Range("E36").Value = "  +5.70%  "

$ws.Range("E37").Value = "  -2.82%  "

$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").Value = "'2.63"
$ws.Range("E39").Value = "  +2.03%  "

$ws.Range("E40").Value = "  -1.64%  "

$ws.Range("D41").Value = "'0.516"
$ws.Range("E41").Value = "  -2.78%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("E44").Value = "  -1.31%  "

$ws.Range("D45").Value = "'0.0466"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("D46").Value = "'5.38"
$ws.Range("E46").Value = "  -4.08%  "

$ws.Range("D47").Value = "'62.34"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").Value = "'0.919"

$ws.Range("D49").Value = "1.710.77"
$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("D51").Value = "'85.28"
$ws.Range("E51").Value = "  -0.86%  "
